# Hma.Calc.xlsx: rename "Index" column to "i" and switch it from a
# 1-based row index (1..501) to a 0-based index (0..500).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell A1: "Index" -> "i"
# (the "testdata" table's column name follows the header cell automatically)
$ws.Range("A1").Value = "i"

# Shift every data row's index down by one: 1..501 -> 0..500
for ($r = 2; $r -le 503; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 2
}

# Column A now only needs to fit up to 3 digits, so it narrows slightly.
$ws.Columns.Item(1).ColumnWidth = 3.17
